# Update latest output (run 142)

$wb = $excel.ActiveWorkbook

# ---- Schedule sheet ----
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E3").Value = -202.74691125
$schedule.Range("F3").Value = -6.704593625992064
$schedule.Range("E4").Value = 501.08655675
$schedule.Range("F4").Value = 33.14064528769841

# ---- Detailed sheet ----
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B15").Value = 77.94
$detailed.Range("B16").Value = 35.88

$detailed.Range("B17").Value = 56.49932
$detailed.Range("C17").Value = "historical"

$detailed.Range("B18").Value = 3.32856
$detailed.Range("C18").Value = "historical"

$detailed.Range("B19").Value = -0.97745
$detailed.Range("B20").Value = -6.49292
$detailed.Range("B21").Value = -7.02733
$detailed.Range("B22").Value = -8.178520000000001
$detailed.Range("B23").Value = -10.56532
$detailed.Range("B24").Value = -14
$detailed.Range("B25").Value = -14.21514
$detailed.Range("B26").Value = -14
$detailed.Range("B27").Value = -18.76135
$detailed.Range("B28").Value = -18.35215
$detailed.Range("B29").Value = -21.90547
$detailed.Range("B30").Value = -23.5
$detailed.Range("B31").Value = -21.38033
$detailed.Range("B32").Value = -21.91813
$detailed.Range("B33").Value = -10
$detailed.Range("B34").Value = -5.74088
$detailed.Range("B35").Value = -0.89533

$detailed.Range("B37").Value = 47.13993
$detailed.Range("B38").Value = 54.727
$detailed.Range("B39").Value = 64.8901
$detailed.Range("B40").Value = 73.19

$detailed.Range("B44").Value = 57.94958
$detailed.Range("B45").Value = 61.8265
$detailed.Range("B47").Value = 63.95875
